$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Knight+Player")

# --- Blank bold-styled separator row beneath the existing "Scripts pt 2" line ---
$ws.Rows.Item(72).Insert()

# --- New "ONCE ASSETS FIXED:" checklist block ---
$ws.Rows.Item(73).Insert()
$ws.Range("A73").Value = "ONCE ASSETS FIXED:"

# --- Sign off on Nick's "All Sounds" entry (row 69) ---
$ws.Range("C69").Value = "Nick 2:27"

$ws.Rows.Item(74).Insert()
$ws.Range("A74").Value = "Begin working in sandbox"

$ws.Rows.Item(75).Insert()
$ws.Range("A75").Value = "Implement files individually"

$ws.Rows.Item(76).Insert()
$ws.Range("A76").Value = "Sounds:"

$ws.Rows.Item(77).Insert()
$ws.Range("A77").Value = "Sprites:"

$ws.Rows.Item(78).Insert()
$ws.Range("A78").Value = "Scripts:"

# Row 79 stays blank/unused between the two checklists - insert it only to
# pick up the bold row styling for row 80, then clear it out entirely so no
# row 79 is left behind.
$ws.Rows.Item(79).Insert()
$ws.Rows.Item(80).Insert()
$ws.Range("A80").Value = "ONCE ASSETS IMPLEMENTED:"

$ws.Rows.Item(81).Insert()
$ws.Range("A81").Value = "Begin combining sandboxes"

$ws.Rows.Item(79).Clear()

# --- Update the view to reflect the newly-added content at the bottom ---
$excel.ActiveWindow.ScrollRow = 59
$ws.Range("A82").Select() | Out-Null
